$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a simulation-scheme x HKL/pairing grid:
#   row 1   -> running column index (0,1,2,...)
#   row 2   -> column headers (HKL reflection triples, then pairing labels)
#   col A   -> running row index
#   col B   -> simulation-scheme name
#   C:* the data grid (all 1s)
#
# Columns X:AG were a verbatim duplicate of columns N:W (left over cruft) -
# they are removed. Four new simulation rows ("Holden2.5/5/10/15") are
# introduced in place of the old "HexGrid-90degTilt*" rows, and the
# "HexGrid-90degTilt*" data is re-appended as brand-new rows 20-23. Row 2's
# HKL headers (C2:M2) are re-ordered to match the new scheme layout.
# ---------------------------------------------------------------------------

# 1. Drop the stray duplicate columns X:AG (mirrored N:W) entirely.
$ws.Range("X1:AG19").EntireColumn.Delete()

# 2. Append the old "HexGrid-90degTilt*" rows (currently 16:19) as new rows
#    20:23, carrying over their values and formatting unchanged.
$ws.Range("A16:W19").Copy($ws.Range("A20:W23"))

# Fix up the running index in column A for the newly appended rows.
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21

# 3. Re-purpose rows 16:19 (same data, same formatting) as the new
#    "Holden" simulation scheme.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 4. Re-order the HKL reflection headers in row 2 (columns C:M) to match
#    the new scheme's ordering.
$ws.Range("C2").Value = "[5, 1, 1]"
$ws.Range("D2").Value = "[4, 2, 2]"
$ws.Range("E2").Value = "[3, 3, 1]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[1, 1, 1]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 3, 3]"
$ws.Range("J2").Value = "[2, 0, 0]"
$ws.Range("K2").Value = "[2, 2, 0]"
$ws.Range("L2").Value = "[4, 2, 0]"
$ws.Range("M2").Value = "[4, 0, 0]"
